# Auto-generated edit script applying F-column ("想去人数" / want-to-go count) updates
# across all four worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5715
$ws.Range("F7").Value = 8163
$ws.Range("F10").Value = 3884
$ws.Range("F12").Value = 26
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 208
$ws.Range("F15").Value = 140
$ws.Range("F16").Value = 30
$ws.Range("F18").Value = 109
$ws.Range("F20").Value = 622
$ws.Range("F21").Value = 3910
$ws.Range("F24").Value = 5354
$ws.Range("F26").Value = 2122
$ws.Range("F27").Value = 134
$ws.Range("F28").Value = 358
$ws.Range("F29").Value = 7971
$ws.Range("F33").Value = 2204
$ws.Range("F34").Value = 1336
$ws.Range("F35").Value = 1308
$ws.Range("F37").Value = 26
$ws.Range("F38").Value = 272
$ws.Range("F41").Value = 1183
$ws.Range("F42").Value = 1178
$ws.Range("F44").Value = 1341
$ws.Range("F45").Value = 2103
$ws.Range("F46").Value = 134

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 148
$ws.Range("F11").Value = 125
$ws.Range("F20").Value = 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 582
$ws.Range("F3").Value = 758

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 758
$ws.Range("F6").Value = 5715
$ws.Range("F8").Value = 3884
$ws.Range("F10").Value = 26
$ws.Range("F11").Value = 28
$ws.Range("F12").Value = 140
$ws.Range("F13").Value = 30
$ws.Range("F15").Value = 109
$ws.Range("F16").Value = 148
$ws.Range("F18").Value = 622
$ws.Range("F19").Value = 3910
$ws.Range("F23").Value = 5354
$ws.Range("F25").Value = 2122
$ws.Range("F26").Value = 134
$ws.Range("F27").Value = 358
$ws.Range("F28").Value = 7971
$ws.Range("F31").Value = 2204
$ws.Range("F32").Value = 1336
$ws.Range("F33").Value = 1308
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 272
$ws.Range("F38").Value = 1183
$ws.Range("F39").Value = 1178
$ws.Range("F42").Value = 1341
$ws.Range("F44").Value = 2103
$ws.Range("F45").Value = 134
$ws.Range("F48").Value = 15
